# Update the "Training Dashboard" sheet with the new progress date (04-Nov-2025).
# For every data row (3 through 40):
#   - Column H "PERIOD TO EXPIRE" decreases by 1 day
#   - Column I "LAST UPDATE" changes from 03-Nov-2025 to 04-Nov-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

$firstRow = 3
$lastRow = 40

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $hCell = $ws.Cells.Item($r, 8)   # Column H
    $iCell = $ws.Cells.Item($r, 9)   # Column I

    $hVal = $hCell.Value2
    if ($hVal -ne $null) {
        $hCell.Value = $hVal - 1
    }

    $iVal = $iCell.Value2
    if ($iVal -eq "03-Nov-2025") {
        # Force the cell to stay as text so Excel does not silently
        # reinterpret the "dd-mmm-yyyy" looking string as a date value.
        $iCell.NumberFormat = "@"
        $iCell.Value = "04-Nov-2025"
    }
}
